$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 58, shifting existing rows 58-135 down to 59-136.
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the latest weekly price observation.
$ws.Range("A58").Value = 11
$ws.Range("B58").Value = "Vega Monumental Concepción"
$ws.Range("C58").Value = "Bíobío"
$ws.Range("D58").Value = 44483
$ws.Range("D58").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 100114013
$ws.Range("G58").Value = "Zanahoria"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 550
$ws.Range("K58").Value = 8000
$ws.Range("L58").Value = 8500
$ws.Range("M58").Value = 8227
$ws.Range("N58").Value = "`$/saco 20 kilos"
$ws.Range("O58").Value = "Chillán"
$ws.Range("P58").Value = 411
$ws.Range("Q58").Value = 20
$ws.Range("R58").Value = "Hortaliza"
